$wb = $excel.ActiveWorkbook

# "Repayment schedule" is the 3rd sheet (sheet3.xml / rId3)
$ws = $wb.Worksheets.Item(3)

# Insert a new blank column before the existing "Late" column (old column N),
# shifting Late/Waived/Outstanding one column to the right (N->O, O->P, P->Q)
$ws.Columns("N:N").Insert() | Out-Null
$ws.Columns("N:N").ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet/tab (was "Transactions")
$ws.Activate() | Out-Null

# Update the selected cell on the "Repayment schedule" sheet
$ws.Range("K17").Select() | Out-Null
